# Preparation for the 5/14/13 status meeting.
# Update existing action-item rows 79-81 (status changes) and append two new
# action items (rows 82-83) plus one trailing blank row (84), matching the
# caArray/caIntegrator Action Items tracker conventions already used in the
# sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 79 & 80: adopt the "most recently closed-out" banding/border used
# by row 78 (copy formats only, values are changed separately below). ---
foreach ($col in @("A","B","C","D","E")) {
    $ws.Range("$col`78").Copy()
    $ws.Range("$col`79").PasteSpecial(-4122)
}
foreach ($col in @("A","B","C","D","E")) {
    $ws.Range("$col`78").Copy()
    $ws.Range("$col`80").PasteSpecial(-4122)
}

# Status updates: items #78 and #79 are now Complete.
$ws.Range("E79").Value = "Complete"
$ws.Range("E80").Value = "Complete"

# Item #80 (gitHub integration on caIntegrator PRODUCTION) moved to In Progress.
$ws.Range("E81").Value = "In Progress"

# --- New rows 82 & 83: same style/banding as row 81 (which keeps its
# original formatting), plus a trailing blank row 84 in the same style. ---
foreach ($col in @("A","B","C","D","E")) {
    $ws.Range("$col`81").Copy()
    $ws.Range("$col`82").PasteSpecial(-4122)
}
foreach ($col in @("A","B","C","D","E")) {
    $ws.Range("$col`81").Copy()
    $ws.Range("$col`83").PasteSpecial(-4122)
}
foreach ($col in @("A","B","C","D","E")) {
    $ws.Range("$col`81").Copy()
    $ws.Range("$col`84").PasteSpecial(-4122)
}

$ws.Range("A82").Value = 81
$ws.Range("B82").Value = "Confirm that all caArray tiers now have gitHub integration with AntHill Pro configured."
$ws.Range("C82").Value = "Mike Hunter"
$ws.Range("D82").Value = Get-Date -Year 2013 -Month 5 -Day 7
$ws.Range("E82").Value = "In Progress"

$ws.Range("A83").Value = 82
$ws.Range("B83").Value = "Hold meeting with Juli Klemm, Ulli Wagner, and JJ Pan to plan the caArray 2.5.3 release."
$ws.Range("C83").Value = "Mike Hunter"
$ws.Range("D83").Value = Get-Date -Year 2013 -Month 5 -Day 14
$ws.Range("E83").Value = "Assigned"

# Row 84 stays blank (already cleared by the PasteSpecial(-4122) which only
# copied formatting, not values).

# Match the saved selection/viewport from the commit: cursor left on C83.
$ws.Range("C83").Select()
